# Updated LoadExcelData and added TimeSeries script
# Appends additional traffic-light sample rows (Time=3..9) to Sheet1,
# matching the style of the existing data rows, and updates the active
# selection to match where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: Time (s), Car, Bus, Emergency
$newRows = @(
    @(3, 15, 3, 0),
    @(4, 17, 5, 0),
    @(5, 18, 3, 2),
    @(6, 15, 3, 0),
    @(7, 14, 4, 0),
    @(8, 13, 2, 1),
    @(9, 12, 1, 0)
)

$startRow = 5
$endRow = $startRow + $newRows.Count - 1
$targetAddress = "A{0}:D{1}" -f $startRow, $endRow

# Give the new rows the same formatting as the existing data rows by
# copying row 2's format down across A5:D11
$ws.Range("A2:D2").Copy()
$ws.Range($targetAddress).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Update the selection to mirror where the editor left the cursor (E9)
$ws.Range("E9").Select()
